# Auto-generated Excel COM-interop edit script
# Applies numeric value updates to columns H-N across multiple sheets
# as described by the upstream OOXML diff (Sargatanas_Profits.xlsx).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 83.5
$ws.Range("I2").Value = 66.666664
$ws.Range("K2").Value = 66.666664
$ws.Range("M2").Value = 46.333336
$ws.Range("H18").Value = 6468.125
$ws.Range("J18").Value = 5449.5
$ws.Range("L18").Value = 5449.5
$ws.Range("N18").Value = -6017.5
$ws.Range("H33").Value = 708.1818
$ws.Range("I33").Value = 833.9375
$ws.Range("J33").Value = 372.83334
$ws.Range("K33").Value = 833.9375
$ws.Range("L33").Value = 372.83334
$ws.Range("M33").Value = -604.9375
$ws.Range("N33").Value = -830.83334
$ws.Range("H40").Value = 3500
$ws.Range("I40").Value = 3500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -3325
$ws.Range("H53").Value = 3027.0908
$ws.Range("I53").Value = 2517.5386
$ws.Range("J53").Value = 3763.111
$ws.Range("K53").Value = 2517.5386
$ws.Range("L53").Value = 3763.111
$ws.Range("M53").Value = -1880.5386
$ws.Range("N53").Value = -5037.111
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = $null
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = $null
$ws.Range("N67").Value = 0
$ws.Range("H86").Value = 83335784
$ws.Range("I86").Value = 103176450
$ws.Range("J86").Value = 13893455
$ws.Range("K86").Value = 103176450
$ws.Range("L86").Value = 13893455
$ws.Range("M86").Value = -103175327
$ws.Range("N86").Value = -13895701
$ws.Range("H89").Value = 83335784
$ws.Range("I89").Value = 103176450
$ws.Range("J89").Value = 13893455
$ws.Range("K89").Value = 515882250
$ws.Range("L89").Value = 69467275
$ws.Range("M89").Value = -515876634
$ws.Range("N89").Value = -69478507
$ws.Range("H133").Value = 129592
$ws.Range("J133").Value = 129592
$ws.Range("L133").Value = 129592
$ws.Range("N133").Value = -139712
$ws.Range("H136").Value = 53593.332
$ws.Range("J136").Value = 60780
$ws.Range("L136").Value = 60780
$ws.Range("N136").Value = -70980
$ws.Range("H137").Value = 2566.524
$ws.Range("I137").Value = 2571.3845
$ws.Range("J137").Value = 2558.625
$ws.Range("K137").Value = 7714.1535
$ws.Range("L137").Value = 7675.875
$ws.Range("M137").Value = -5164.1535
$ws.Range("N137").Value = -12775.875
$ws.Range("H138").Value = 6877
$ws.Range("I138").Value = 2854.2222
$ws.Range("J138").Value = 8522.682000000001
$ws.Range("K138").Value = 8562.6666
$ws.Range("L138").Value = 25568.046
$ws.Range("M138").Value = -3422.6666
$ws.Range("N138").Value = -35848.046

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 100377
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 100377
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = $null
$ws.Range("M43").Value = 100377
$ws.Range("N43").Value = -101003
$ws.Range("H61").Value = 5411.4316
$ws.Range("I61").Value = 2544.9312
$ws.Range("K61").Value = 2544.9312
$ws.Range("M61").Value = -2332.9312
$ws.Range("H114").Value = 55652.668
$ws.Range("J114").Value = 55652.668
$ws.Range("L114").Value = 55652.668
$ws.Range("N114").Value = -64330.668
$ws.Range("H119").Value = 67319
$ws.Range("J119").Value = 67319
$ws.Range("L119").Value = 67319
$ws.Range("N119").Value = -76995
$ws.Range("H132").Value = 6854
$ws.Range("J132").Value = 9079
$ws.Range("L132").Value = 27237
$ws.Range("N132").Value = -32297
$ws.Range("H136").Value = 5411.4316
$ws.Range("I136").Value = 2544.9312
$ws.Range("K136").Value = 7634.7936
$ws.Range("M136").Value = -5084.7936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 41669710
$ws.Range("I107").Value = 51138720
$ws.Range("K107").Value = 51138720
$ws.Range("M107").Value = -51136800
$ws.Range("H134").Value = 4802.5
$ws.Range("I134").Value = 1856.8485
$ws.Range("J134").Value = 9028.869000000001
$ws.Range("K134").Value = 5570.5455
$ws.Range("L134").Value = 27086.607
$ws.Range("M134").Value = -3035.5455
$ws.Range("N134").Value = -32156.607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8471.451999999999
$ws.Range("J31").Value = 14114.952
$ws.Range("L31").Value = 14114.952
$ws.Range("N31").Value = -14704.952
$ws.Range("H34").Value = 8471.451999999999
$ws.Range("J34").Value = 14114.952
$ws.Range("L34").Value = 14114.952
$ws.Range("N34").Value = -14518.952

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 10499.823
$ws.Range("J132").Value = 17999.715
$ws.Range("L132").Value = 161997.435
$ws.Range("N132").Value = -167057.435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5555.6313
$ws.Range("I102").Value = 4970.8
$ws.Range("K102").Value = 4970.8
$ws.Range("M102").Value = -3348.8
$ws.Range("H121").Value = 46192.75
$ws.Range("J121").Value = 46192.75
$ws.Range("L121").Value = 46192.75
$ws.Range("N121").Value = -49686.75
$ws.Range("H132").Value = 13262.866
$ws.Range("I132").Value = 5176.8335
$ws.Range("J132").Value = 18653.555
$ws.Range("K132").Value = 15530.5005
$ws.Range("L132").Value = 55960.665
$ws.Range("M132").Value = -13000.5005
$ws.Range("N132").Value = -61020.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2666.1428
$ws.Range("I22").Value = 1267
$ws.Range("J22").Value = 3047.7273
$ws.Range("K22").Value = 1267
$ws.Range("L22").Value = 3047.7273
$ws.Range("M22").Value = -972
$ws.Range("N22").Value = -3637.7273
$ws.Range("H27").Value = 2666.1428
$ws.Range("I27").Value = 1267
$ws.Range("J27").Value = 3047.7273
$ws.Range("K27").Value = 1267
$ws.Range("L27").Value = 3047.7273
$ws.Range("M27").Value = -1160
$ws.Range("N27").Value = -3261.7273
$ws.Range("H46").Value = 11113190
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = $null
$ws.Range("H55").Value = 55556132
$ws.Range("H80").Value = 74999
$ws.Range("J80").Value = 74999
$ws.Range("L80").Value = 74999
$ws.Range("N80").Value = -77245
$ws.Range("H83").Value = 74999
$ws.Range("J83").Value = 74999
$ws.Range("L83").Value = 224997
$ws.Range("N83").Value = -236229
$ws.Range("H93").Value = 4800.1113
$ws.Range("I93").Value = 6099.2
$ws.Range("J93").Value = 3176.25
$ws.Range("K93").Value = 6099.2
$ws.Range("L93").Value = 3176.25
$ws.Range("M93").Value = -4851.2
$ws.Range("N93").Value = -5672.25
$ws.Range("H94").Value = 46969
$ws.Range("J94").Value = 46969
$ws.Range("L94").Value = 46969
$ws.Range("N94").Value = -48321
$ws.Range("H119").Value = 56302.332
$ws.Range("J119").Value = 56302.332
$ws.Range("L119").Value = 56302.332
$ws.Range("N119").Value = -65978.33199999999
$ws.Range("H136").Value = 7571.1025
$ws.Range("I136").Value = 1711.5
$ws.Range("K136").Value = 5134.5
$ws.Range("M136").Value = -2584.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14454.546
$ws.Range("I54").Value = 15000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14480
$ws.Range("H113").Value = 1499
$ws.Range("I113").Value = 972.125
$ws.Range("J113").Value = 2025.875
$ws.Range("K113").Value = 2916.375
$ws.Range("L113").Value = 6077.625
$ws.Range("M113").Value = -746.375
$ws.Range("N113").Value = -10417.625
$ws.Range("H122").Value = 338442.34
$ws.Range("I122").Value = 575486.9
$ws.Range("J122").Value = 6580
$ws.Range("K122").Value = 1726460.7
$ws.Range("L122").Value = 19740
$ws.Range("M122").Value = -1724010.7
$ws.Range("N122").Value = -24640
$ws.Range("H135").Value = 72357.5
$ws.Range("J135").Value = 72357.5
$ws.Range("L135").Value = 72357.5
$ws.Range("N135").Value = -82497.5
